# AchAuthLogHistory.xlsx edit
# Updates the "CdCode:" reference-note cells on the "DBD" sheet so that the
# separator character used between a code and its meaning is changed from a
# mix of "-", "." or nothing to a consistent ":" and the prefix separator is
# changed from "CdCode:" to "CdCode." (e.g. G19, G20, G27, G31, G32).
# Finally the active selection is moved to G33, matching the author's last
# selected cell when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# G19 - CdCode:AuthStatus -> CdCode.AuthStatus (dash "-" separators -> ":")
$ws.Range("G19").Value = "CdCode.AuthStatus`n空白:未授權`n0:成功授權/取消授權`n1:印鑑不符`n2:無此帳號`n3:委繳戶統一編號不符`n4:已核印成功在案`n5:原交易不存在`n6:電子資料與授權書內容不符`n7:帳戶已結清`n8:印鑑不清`n9:其他`nA:未收到授權書`nB:用戶號碼錯誤`nC:靜止戶`nD:未收到聲明書`nE:授權書資料不全`nF:警示戶`nG:本帳戶不適用授權扣繳`nH:已於他行授權扣款`nI:該用戶已死亡`nZ:未交易或匯入失敗資料"

# G20 - CdCode:AchAuthCode -> CdCode.AchAuthCode
$ws.Range("G20").Value = "CdCode.AchAuthCode`nA:紙本新增`nO:舊檔轉換"

# G27 - CdCode:RelationCode -> CdCode.RelationCode (add ":" between code and meaning)
$ws.Range("G27").Value = "CdCode.RelationCode`n00:本人`n01:夫`n02:妻`n03:父`n04:母`n05:子`n06:女`n07:兄`n08:弟`n09:姊`n10:妹`n11:姪子`n99:其他"

# G31 - CdCode:Sex -> CdCode.Sex
$ws.Range("G31").Value = "CdCode.Sex"

# G32 - CdCode:AmlCheckItem -> CdCode.AmlCheckItem (dot "." separators -> ":")
$ws.Range("G32").Value = "CdCode.AmlCheckItem`n0:非可疑名單/已完成名單確認`n1:需審查/確認`n2:為凍結名單/未確定名單"

# Move the active selection to G33, as it was when the author last saved.
$ws.Range("G33").Select()
